$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full target roster table (header row unchanged; rows 2-19 reordered/updated
# so each player lines up with their correct position & team).
$data = @(
    @("Russell Westbrook", "PG",       "Denver Nuggets"),
    @("Jalen Green",       "PG,SG",    "Houston Rockets"),
    @("Chris Paul",        "PG",       "San Antonio Spurs"),
    @("Ayo Dosunmu",       "PG,SG,SF", "Chicago Bulls"),
    @("Jaylen Brown",      "SG,SF",    "Boston Celtics"),
    @("Pascal Siakam",     "SF,PF,C",  "Indiana Pacers"),
    @("Deni Avdija",       "SF,PF",    "Portland Trail Blazers"),
    @("Jerami Grant",      "SF,PF",    "Portland Trail Blazers"),
    @("Alexandre Sarr",    "PF,C",     "Washington Wizards"),
    @("Nikola Jokic",      "C",        "Denver Nuggets"),
    @("Dejounte Murray",   "PG,SG",    "New Orleans Pelicans"),
    @("Malcolm Brogdon",   "PG,SG",    "Washington Wizards"),
    @("Jalen Suggs",       "PG,SG",    "Orlando Magic"),
    @("Dru Smith",         "PG,SG",    "Miami Heat"),
    @("Rudy Gobert",       "C",        "Minnesota Timberwolves"),
    @("Paolo Banchero",    "SF,PF",    "Orlando Magic"),
    @("Chet Holmgren",     "PF,C",     "Oklahoma City Thunder"),
    @("Jakob Poeltl",      "C",        "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
